$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for data rows 2 through 79
# from serial date 45172 to serial date 45175.
$ws.Range("C2:C79").Value = 45175
